$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, matching the style of the other header cells (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-33
$values = @{
    2  = @(7, 8)
    3  = @(7, 7)
    4  = @(9, 9)
    5  = @(4, 4)
    6  = @(7, 7)
    7  = @(6, 7)
    8  = @(8, 8)
    9  = @(7, 7)
    10 = @(7, 7)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(6, 6)
    14 = @(8, 8)
    15 = @(4, 5)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(7, 7)
    19 = @(7, 7)
    20 = @(5, 6)
    21 = @(7, 7)
    22 = @(8, 8)
    23 = @(7, 7)
    24 = @(9, 9)
    25 = @(8, 9)
    26 = @(7, 7)
    27 = @(8, 9)
    28 = @(8, 8)
    29 = @(9, 9)
    30 = @(5, 5)
    31 = @(2, 2)
    32 = @(6, 6)
    33 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
